$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.322.03'
$ws.Range("E2").Value = '  +2.63%  '

$ws.Range("D3").Value = '1.870.66'
$ws.Range("E3").Value = '  +1.21%  '

$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.84'
$ws.Range("E5").Value = '  +2.13%  '

$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4707'
$ws.Range("E7").Value = '  +1.58%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3938'
$ws.Range("E8").Value = '  +2.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.30'

$ws.Range("E10").Value = '  +1.01%  '

$ws.Range("E11").Value = '  +1.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.93'
$ws.Range("E12").Value = '  +2.12%  '

$ws.Range("D13").Value = '1.882.11'
$ws.Range("E13").Value = '  +1.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.007'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.286'
$ws.Range("E15").Value = '  +2.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.22'
$ws.Range("E16").Value = '  +2.81%  '

$ws.Range("E17").Value = '  -0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001043'
$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06594'
$ws.Range("E19").Value = '  -0.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.69'
$ws.Range("E20").Value = '  +3.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("D22").Value = '28.324.81'
$ws.Range("E22").Value = '  +2.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.462'
$ws.Range("E23").Value = '  +1.46%  '

$ws.Range("E24").Value = '  +1.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.291'
$ws.Range("E25").Value = '  -0.56%  '

$ws.Range("D26").Value = '2.105.02'
$ws.Range("E26").Value = '  +1.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.97'
$ws.Range("E27").Value = '  +1.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.84'
$ws.Range("E28").Value = '  +1.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.149'
$ws.Range("E29").Value = '  +2.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.501'
$ws.Range("E30").Value = '  +1.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.42'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9769'
$ws.Range("E32").Value = '  +0.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09503'
$ws.Range("E33").Value = '  +1.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.594'
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.381'
$ws.Range("E35").Value = '  +2.72%  '

$ws.Range("E36").Value = '  +1.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02279'
$ws.Range("E37").Value = '  +2.59%  '

$ws.Range("E38").Value = '  +1.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.457'
$ws.Range("E39").Value = '  +1.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.180'
$ws.Range("E40").Value = '  +0.19%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5984'
$ws.Range("E41").Value = '  +1.52%  '

$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1882'
$ws.Range("E43").Value = '  +1.11%  '

$ws.Range("E44").Value = '  +1.19%  '

$ws.Range("E45").Value = '  +5.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5630'
$ws.Range("E46").Value = '  +0.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.11'
$ws.Range("E47").Value = '  -0.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.969'
$ws.Range("E48").Value = '  +3.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06917'
$ws.Range("E49").Value = '  +3.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.36'
$ws.Range("E50").Value = '  +0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.026'
$ws.Range("E51").Value = '  +13.43%  '
